# Form the consolidated report: fill in the "Absent" column (H) for each
# date row based on the "Real" column (E). A student is considered
# Absent (H = 1) on a date when they were not marked Real (E = 0),
# otherwise they were present (H = 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 3; $r -le $lastRow; $r++) {
    $realCell = $ws.Cells.Item($r, 5)   # column E - Real
    $realValue = $realCell.Value2
    if ($realValue -eq $null -or $realValue -eq "") {
        continue
    }
    $absentCell = $ws.Cells.Item($r, 8)  # column H - Absent
    if ([double]$realValue -eq 0) {
        $absentCell.Value = 1
    } else {
        $absentCell.Value = 0
    }
}
